$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.489914298057556
$ws.Range("B1").Value = 1.781668901443481
$ws.Range("C1").Value = 1.89094865322113
$ws.Range("D1").Value = 2.184774160385132
$ws.Range("E1").Value = 2.750805139541626
